$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = -0.3102236102782091
$ws.Range("J2").Value = 0.003511609830582043
$ws.Range("S2").Value = -0.03655835111570035
$ws.Range("U2").Value = 0.03655835111570035
$ws.Range("C3").Value = 59.81720824442149
$ws.Range("G3").Value = 0.06528219984543301
$ws.Range("J3").Value = 0.1723260399241088
$ws.Range("S3").Value = -1.626999868594065
$ws.Range("U3").Value = 1.626999868594065
$ws.Range("C4").Value = 51.68220890145116
$ws.Range("G4").Value = 0.3483097607009766
$ws.Range("J4").Value = 0.3102199171932781
$ws.Range("S4").Value = -2.538001899746587
$ws.Range("U4").Value = 2.538001899746587
$ws.Range("C5").Value = 38.99219940271823
$ws.Range("G5").Value = -0.280623952189323
$ws.Range("J5").Value = 0.004383993456098986
$ws.Range("S5").Value = -0.04885777159885044
$ws.Range("U5").Value = 0.04885777159885044
$ws.Range("C6").Value = 38.74791054472398
$ws.Range("G6").Value = -0.2443879286041725
$ws.Range("J6").Value = 0.003097111546940342
$ws.Range("S6").Value = -0.039711619986608
$ws.Range("U6").Value = 0.039711619986608
$ws.Range("C7").Value = 38.54935244479094
$ws.Range("G7").Value = -0.05167717441882491
$ws.Range("J7").Value = 0.1519658903576547
$ws.Range("S7").Value = -1.467430259851106
$ws.Range("U7").Value = 1.468050040070886
$ws.Range("C8").Value = 31.21220114553541
$ws.Range("G8").Value = 0.4280462598895027
$ws.Range("J8").Value = 0.4519502217217179
$ws.Range("S8").Value = -3.525435370022786
$ws.Range("U8").Value = 3.721107897495313
$ws.Range("C9").Value = 13.58502429542148
$ws.Range("D9").Value = 6.504631141470997
$ws.Range("F9").Value = 0.3252315570735498
$ws.Range("G9").Value = -0.8765221351827451
$ws.Range("J9").Value = 0.06635029916973001
$ws.Range("O9").Value = 0.6114993896827655
$ws.Range("P9").Value = 5.893131751788231
$ws.Range("S9").Value = -0
$ws.Range("T9").Value = -0.7945676432842671
$ws.Range("U9").Value = 0.7945676432842671
$ws.Range("C10").Value = 46.10818000277646
$ws.Range("D10").Value = 6.460424613588851
$ws.Range("F10").Value = 0.3230212306794426
$ws.Range("G10").Value = -0.5089531785430129
$ws.Range("J10").Value = 0.2391949748607684
$ws.Range("O10").Value = 0.05857885668394358
$ws.Range("P10").Value = 6.401845756904907
$ws.Range("S10").Value = -0
$ws.Range("T10").Value = -2.862913451008365
$ws.Range("U10").Value = 2.862913451008365
$ws.Range("C11").Value = 78.41030307072072
$ws.Range("G11").Value = 0.5064483076277337
$ws.Range("J11").Value = 0.4540693482500791
$ws.Range("S11").Value = -0.1076284392677893
$ws.Range("U11").Value = 4.220346021685372
$ws.Range("C12").Value = 77.87216087438178
$ws.Range("G12").Value = 0.6534817408827137
$ws.Range("J12").Value = 0.5790382908318221
$ws.Range("S12").Value = -0.5198184548984184
$ws.Range("U12").Value = 5.194660213140176
$ws.Range("C13").Value = 75.2730685998897
$ws.Range("G13").Value = 1.011564887143696
$ws.Range("J13").Value = 0.846816606577889
$ws.Range("S13").Value = -1.044740180188295
$ws.Range("U13").Value = 6.047628092276207
$ws.Range("C14").Value = 70.04936769894822
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 1.015491394397412
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0.8530509685025724
$ws.Range("R14").Value = 0
$ws.Range("S14").Value = -1.675214413169249
$ws.Range("U14").Value = 6.373994632949469
$ws.Range("C15").Value = 61.67329563310198
$ws.Range("G15").Value = 0.6253477224934187
$ws.Range("J15").Value = 0.5649835139638191
$ws.Range("S15").Value = -1.810415586881421
$ws.Range("U15").Value = 6.038127674793508
$ws.Range("C16").Value = 52.62121769869488
$ws.Range("D16").Value = 7.475756460261024
$ws.Range("F16").Value = 0.3737878230130512
$ws.Range("G16").Value = -0.5901451629441445
$ws.Range("J16").Value = 0.2280097693561362
$ws.Range("O16").Value = 0.01772606941523147
$ws.Range("P16").Value = 7.458030390845793
$ws.Range("T16").Value = -3.217523381134219
$ws.Range("U16").Value = 3.217523381134219
$ws.Range("C17").Value = 90
$ws.Range("G17").Value = 0.8796003027776238
$ws.Range("J17").Value = 0.7423994118356929
$ws.Range("S17").Value = -2.863972435302911
$ws.Range("U17").Value = 4.460371336401812
$ws.Range("C18").Value = 75.68013782348544
$ws.Range("D18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0.6398642749122446
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0.5473047485497745
$ws.Range("R18").Value = 0
$ws.Range("S18").Value = -3.751872418793567
$ws.Range("U18").Value = 4.244270220991369
$ws.Range("C19").Value = 56.9207757295176
$ws.Range("G19").Value = -0.4849647841108055
$ws.Range("J19").Value = 0.05474283345783239
$ws.Range("S19").Value = -0.6010960640895577
$ws.Range("U19").Value = 0.6358103498038434
$ws.Range("C20").Value = 53.91529540906982
$ws.Range("G20").Value = -0.3314358803587044
$ws.Range("J20").Value = 0.169752435505516
$ws.Range("S20").Value = -1.705576792603539
$ws.Range("U20").Value = 1.705576792603539
$ws.Range("C21").Value = 45.38741144605212
$ws.Range("G21").Value = -0.2909234520802559
$ws.Range("J21").Value = 0.1410894779829487
$ws.Range("S21").Value = -1.507345480961397
$ws.Range("U21").Value = 1.507345480961397
$ws.Range("C22").Value = 37.85068404124514
$ws.Range("G22").Value = -0.1879398448635849
$ws.Range("J22").Value = 0.189500974316387
$ws.Range("S22").Value = -1.893084013986703
$ws.Range("U22").Value = 1.893084013986703
$ws.Range("C23").Value = 28.38526397131162
$ws.Range("G23").Value = 0.3524813871803485
$ws.Range("J23").Value = 0.4492469405581101
$ws.Range("S23").Value = -3.547049023958224
$ws.Range("U23").Value = 3.547049023958224
$ws.Range("C24").Value = 10.65001885152051
$ws.Range("D24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = -0.3513378968690586
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0.007099609852110451
$ws.Range("R24").Value = 0
$ws.Range("S24").Value = -0.08339877708954342
$ws.Range("U24").Value = 0.08339877708954342
$ws.Range("C25").Value = 10.23302496607279
